$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi following Dr Hou advice:
# Rows 2-4 (sending cluster ECs) are updated in place, and rows 5-10 are
# added so the sheet holds the full 3x3 cross of sending/target clusters
# (ECs, FAPs, sCs) instead of only the three "target = ECs" rows.
# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vcan"
$ws.Cells.Item(2, 3).Value = "Selp"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 7.172092999999999
$ws.Cells.Item(2, 8).Value = 21.516279
$ws.Cells.Item(2, 9).Value = 0.073573870768057
$ws.Cells.Item(2, 10).Value = 0.07357387076805699
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 26.450162
$ws.Cells.Item(2, 14).Value = 79.350486
$ws.Cells.Item(2, 15).Value = 0.9821351879331711
$ws.Cells.Item(2, 16).Value = 0.9821351879331711
$ws.Cells.Item(2, 17).Value = 189.703021729066
$ws.Cells.Item(2, 18).Value = 1707.327195561594
$ws.Cells.Item(2, 19).Value = 0.07225948739375651
$ws.Cells.Item(2, 20).Value = 0.0722594873937565
# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vcan"
$ws.Cells.Item(3, 3).Value = "Selp"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 7.172092999999999
$ws.Cells.Item(3, 8).Value = 21.516279
$ws.Cells.Item(3, 9).Value = 0.073573870768057
$ws.Cells.Item(3, 10).Value = 0.07357387076805699
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.321934
$ws.Cells.Item(3, 14).Value = 0.965802
$ws.Cells.Item(3, 15).Value = 0.01195390446349922
$ws.Cells.Item(3, 16).Value = 0.01195390446349922
$ws.Cells.Item(3, 17).Value = 2.308940587862
$ws.Cells.Item(3, 18).Value = 20.780465290758
$ws.Cells.Item(3, 19).Value = 0.0008794950221711915
$ws.Cells.Item(3, 20).Value = 0.0008794950221711915
# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vcan"
$ws.Cells.Item(4, 3).Value = "Selp"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 7.172092999999999
$ws.Cells.Item(4, 8).Value = 21.516279
$ws.Cells.Item(4, 9).Value = 0.073573870768057
$ws.Cells.Item(4, 10).Value = 0.07357387076805699
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1591883333333333
$ws.Cells.Item(4, 14).Value = 0.477565
$ws.Cells.Item(4, 15).Value = 0.005910907603329674
$ws.Cells.Item(4, 16).Value = 0.005910907603329674
$ws.Cells.Item(4, 17).Value = 1.141713531181667
$ws.Cells.Item(4, 18).Value = 10.275421780635
$ws.Cells.Item(4, 19).Value = 0.000434888352129303
$ws.Cells.Item(4, 20).Value = 0.000434888352129303
# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Vcan"
$ws.Cells.Item(5, 3).Value = "Selp"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 89.72947699999999
$ws.Cells.Item(5, 8).Value = 269.188431
$ws.Cells.Item(5, 9).Value = 0.9204767624852804
$ws.Cells.Item(5, 10).Value = 0.9204767624852804
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 26.450162
$ws.Cells.Item(5, 14).Value = 79.350486
$ws.Cells.Item(5, 15).Value = 0.9821351879331711
$ws.Cells.Item(5, 16).Value = 0.9821351879331711
$ws.Cells.Item(5, 17).Value = 2373.359202825274
$ws.Cells.Item(5, 18).Value = 21360.23282542747
$ws.Cells.Item(5, 19).Value = 0.9040326181115979
$ws.Cells.Item(5, 20).Value = 0.9040326181115979
# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Vcan"
$ws.Cells.Item(6, 3).Value = "Selp"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 89.72947699999999
$ws.Cells.Item(6, 8).Value = 269.188431
$ws.Cells.Item(6, 9).Value = 0.9204767624852804
$ws.Cells.Item(6, 10).Value = 0.9204767624852804
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.321934
$ws.Cells.Item(6, 14).Value = 0.965802
$ws.Cells.Item(6, 15).Value = 0.01195390446349922
$ws.Cells.Item(6, 16).Value = 0.01195390446349922
$ws.Cells.Item(6, 17).Value = 28.88696944851799
$ws.Cells.Item(6, 18).Value = 259.982725036662
$ws.Cells.Item(6, 19).Value = 0.01100329127962011
$ws.Cells.Item(6, 20).Value = 0.01100329127962011
# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Vcan"
$ws.Cells.Item(7, 3).Value = "Selp"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 89.72947699999999
$ws.Cells.Item(7, 8).Value = 269.188431
$ws.Cells.Item(7, 9).Value = 0.9204767624852804
$ws.Cells.Item(7, 10).Value = 0.9204767624852804
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.1591883333333333
$ws.Cells.Item(7, 14).Value = 0.477565
$ws.Cells.Item(7, 15).Value = 0.005910907603329674
$ws.Cells.Item(7, 16).Value = 0.005910907603329674
$ws.Cells.Item(7, 17).Value = 14.28388589450167
$ws.Cells.Item(7, 18).Value = 128.554973050515
$ws.Cells.Item(7, 19).Value = 0.005440853094062527
$ws.Cells.Item(7, 20).Value = 0.005440853094062527
# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Vcan"
$ws.Cells.Item(8, 3).Value = "Selp"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.5799533333333334
$ws.Cells.Item(8, 8).Value = 1.73986
$ws.Cells.Item(8, 9).Value = 0.005949366746662454
$ws.Cells.Item(8, 10).Value = 0.005949366746662453
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 26.450162
$ws.Cells.Item(8, 14).Value = 79.350486
$ws.Cells.Item(8, 15).Value = 0.9821351879331711
$ws.Cells.Item(8, 16).Value = 0.9821351879331711
$ws.Cells.Item(8, 17).Value = 15.33985961910667
$ws.Cells.Item(8, 18).Value = 138.05873657196
$ws.Cells.Item(8, 19).Value = 0.005843082427816689
$ws.Cells.Item(8, 20).Value = 0.005843082427816688
# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Vcan"
$ws.Cells.Item(9, 3).Value = "Selp"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.5799533333333334
$ws.Cells.Item(9, 8).Value = 1.73986
$ws.Cells.Item(9, 9).Value = 0.005949366746662454
$ws.Cells.Item(9, 10).Value = 0.005949366746662453
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.321934
$ws.Cells.Item(9, 14).Value = 0.965802
$ws.Cells.Item(9, 15).Value = 0.01195390446349922
$ws.Cells.Item(9, 16).Value = 0.01195390446349922
$ws.Cells.Item(9, 17).Value = 0.1867066964133334
$ws.Cells.Item(9, 18).Value = 1.68036026772
$ws.Cells.Item(9, 19).Value = 0.00007111816170792215
$ws.Cells.Item(9, 20).Value = 0.00007111816170792215
# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Vcan"
$ws.Cells.Item(10, 3).Value = "Selp"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.5799533333333334
$ws.Cells.Item(10, 8).Value = 1.73986
$ws.Cells.Item(10, 9).Value = 0.005949366746662454
$ws.Cells.Item(10, 10).Value = 0.005949366746662453
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.1591883333333333
$ws.Cells.Item(10, 14).Value = 0.477565
$ws.Cells.Item(10, 15).Value = 0.005910907603329674
$ws.Cells.Item(10, 16).Value = 0.005910907603329674
$ws.Cells.Item(10, 17).Value = 0.09232180454444447
$ws.Cells.Item(10, 18).Value = 0.8308962409000001
$ws.Cells.Item(10, 19).Value = 0.00003516615713784383
$ws.Cells.Item(10, 20).Value = 0.00003516615713784382
